$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-33 down to 30-34
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record
$ws.Cells.Item(29, 1).Value = 7
$ws.Cells.Item(29, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(29, 3).Value = "Ñuble"
$ws.Cells.Item(29, 4).Value = 44644
$ws.Cells.Item(29, 5).Value = 16
$ws.Cells.Item(29, 6).Value = 100112040
$ws.Cells.Item(29, 7).Value = "Cilantro"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 550
$ws.Cells.Item(29, 12).Value = 600
$ws.Cells.Item(29, 13).Value = 575
$ws.Cells.Item(29, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(29, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(29, 16).Value = 575
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = "Hortaliza"

# Match the date-time number format used by the other rows in column D
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(30, 4).NumberFormat
